# Auto-generated data refresh for the Leves profit tables.
# For each affected row, update the market-price / profit columns (H:N)
# to the latest cached values; cells that no longer have a value are cleared
# instead of being zero-filled, matching the upstream market-data behaviour.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 4100
$ws.Range("I2").Value = 4100
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4100
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -3987
$ws.Range("N2").ClearContents()

# Row 17
$ws.Range("H17").Value = 1812.5
$ws.Range("J17").Value = 1812.5
$ws.Range("L17").Value = 5437.5
$ws.Range("N17").Value = -5773.5

# Row 58
$ws.Range("H58").Value = 7530.706
$ws.Range("I58").Value = 108.75
$ws.Range("K58").Value = 326.25
$ws.Range("M58").Value = -176.25

# Row 86
$ws.Range("H86").Value = 5194.3335
$ws.Range("J86").Value = 5964.143
$ws.Range("L86").Value = 5964.143
$ws.Range("N86").Value = -8210.143

# Row 89
$ws.Range("H89").Value = 5194.3335
$ws.Range("J89").Value = 5964.143
$ws.Range("L89").Value = 29820.715
$ws.Range("N89").Value = -41052.715

# Row 100
$ws.Range("H100").Value = 992.7692
$ws.Range("I100").Value = 906.5
$ws.Range("J100").Value = 1130.8
$ws.Range("K100").Value = 906.5
$ws.Range("L100").Value = 1130.8
$ws.Range("M100").Value = -365.5
$ws.Range("N100").Value = -2212.8


$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 25934998
$ws.Range("I61").Value = 21828944
$ws.Range("J61").Value = 41674868
$ws.Range("K61").Value = 21828944
$ws.Range("L61").Value = 41674868
$ws.Range("M61").Value = -21828732
$ws.Range("N61").Value = -41675292

# Row 74
$ws.Range("H74").Value = 8342359
$ws.Range("I74").Value = 13161111
$ws.Range("K74").Value = 13161111
$ws.Range("M74").Value = -13160237

# Row 77
$ws.Range("H77").Value = 8342359
$ws.Range("I77").Value = 13161111
$ws.Range("K77").Value = 65805555
$ws.Range("M77").Value = -65801187

# Row 97
$ws.Range("H97").Value = 1086.3611
$ws.Range("I97").Value = 1054.9062
$ws.Range("J97").Value = 1338
$ws.Range("K97").Value = 1054.9062
$ws.Range("L97").Value = 1338
$ws.Range("M97").Value = -558.9061999999999
$ws.Range("N97").Value = -2330

# Row 132
$ws.Range("H132").Value = 13340202
$ws.Range("I132").Value = 19611374
$ws.Range("J132").Value = 13961.625
$ws.Range("K132").Value = 58834122
$ws.Range("L132").Value = 41884.875
$ws.Range("M132").Value = -58831592
$ws.Range("N132").Value = -46944.875

# Row 136
$ws.Range("H136").Value = 25934998
$ws.Range("I136").Value = 21828944
$ws.Range("J136").Value = 41674868
$ws.Range("K136").Value = 65486832
$ws.Range("L136").Value = 125024604
$ws.Range("M136").Value = -65484282
$ws.Range("N136").Value = -125029704


$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2724.75
$ws.Range("I94").Value = 599.3333
$ws.Range("J94").Value = 4000
$ws.Range("K94").Value = 599.3333
$ws.Range("L94").Value = 4000
$ws.Range("M94").Value = -148.3333
$ws.Range("N94").Value = -4902

# Row 134
$ws.Range("H134").Value = 231905.86
$ws.Range("I134").Value = 1761.4193
$ws.Range("J134").Value = 651581
$ws.Range("K134").Value = 5284.257900000001
$ws.Range("L134").Value = 1954743
$ws.Range("M134").Value = -2749.257900000001
$ws.Range("N134").Value = -1959813


$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 4923.095
$ws.Range("I7").Value = 134.66667
$ws.Range("J7").Value = 16894.166
$ws.Range("K7").Value = 134.66667
$ws.Range("L7").Value = 16894.166
$ws.Range("M7").Value = -21.66667000000001
$ws.Range("N7").Value = -17120.166

# Row 23
$ws.Range("H23").Value = 6000
$ws.Range("J23").Value = 6000
$ws.Range("L23").Value = 6000
$ws.Range("N23").Value = -6480

# Row 27
$ws.Range("H27").Value = 6000
$ws.Range("J27").Value = 6000
$ws.Range("L27").Value = 6000
$ws.Range("N27").Value = -6384

# Row 31
$ws.Range("H31").Value = 1233687.5
$ws.Range("I31").Value = 1916.6666
$ws.Range("K31").Value = 1916.6666
$ws.Range("M31").Value = -1621.6666

# Row 34
$ws.Range("H34").Value = 1233687.5
$ws.Range("I34").Value = 1916.6666
$ws.Range("K34").Value = 1916.6666
$ws.Range("M34").Value = -1714.6666

# Row 88
$ws.Range("H88").Value = 99000
$ws.Range("J88").Value = 99000
$ws.Range("L88").Value = 99000
$ws.Range("N88").Value = -99812

# Row 91
$ws.Range("H91").Value = 99000
$ws.Range("J91").Value = 99000
$ws.Range("L91").Value = 99000
$ws.Range("N91").Value = -101808

# Row 105
$ws.Range("H105").Value = 2077
$ws.Range("I105").Value = 2094.875
$ws.Range("K105").Value = 2094.875
$ws.Range("M105").Value = -347.875

# Row 107
$ws.Range("H107").Value = 418.35
$ws.Range("I107").Value = 414.05264
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 414.05264
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1505.94736
$ws.Range("N107").Value = -4340

# Row 120
$ws.Range("H120").Value = 27997
$ws.Range("I120").Value = 19995
$ws.Range("K120").Value = 19995
$ws.Range("M120").Value = -16366

# Row 132
$ws.Range("H132").Value = 10413
$ws.Range("I132").Value = 2212.6667
$ws.Range("K132").Value = 6638.000100000001
$ws.Range("M132").Value = -4108.000100000001

# Row 134
$ws.Range("H134").Value = 5420.5835
$ws.Range("I134").Value = 2425
$ws.Range("J134").Value = 8416.166999999999
$ws.Range("K134").Value = 7275
$ws.Range("L134").Value = 25248.501
$ws.Range("M134").Value = -4740
$ws.Range("N134").Value = -30318.501


$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1873.238
$ws.Range("I5").Value = 1471
$ws.Range("J5").Value = 2034.1333
$ws.Range("K5").Value = 4413
$ws.Range("L5").Value = 6102.3999
$ws.Range("M5").Value = -4301
$ws.Range("N5").Value = -6326.3999

# Row 40
$ws.Range("H40").Value = 74.46154
$ws.Range("I40").Value = 37.22222
$ws.Range("J40").Value = 158.25
$ws.Range("K40").Value = 148.88888
$ws.Range("L40").Value = 633
$ws.Range("M40").Value = -79.88888
$ws.Range("N40").Value = -771

# Row 60
$ws.Range("H60").Value = 1530.125
$ws.Range("I60").Value = 363.2
$ws.Range("K60").Value = 1089.6
$ws.Range("M60").Value = -838.5999999999999

# Row 134
$ws.Range("H134").Value = 13628.066
$ws.Range("I134").Value = 13990.75
$ws.Range("J134").Value = 13213.571
$ws.Range("K134").Value = 41972.25
$ws.Range("L134").Value = 39640.713
$ws.Range("M134").Value = -36902.25
$ws.Range("N134").Value = -49780.713

# Row 135
$ws.Range("H135").Value = 1873.238
$ws.Range("I135").Value = 1471
$ws.Range("J135").Value = 2034.1333
$ws.Range("K135").Value = 13239
$ws.Range("L135").Value = 18307.1997
$ws.Range("M135").Value = -10704
$ws.Range("N135").Value = -23377.1997

# Row 137
$ws.Range("H137").Value = 3771.8572
$ws.Range("J137").Value = 5098.9
$ws.Range("L137").Value = 15296.7
$ws.Range("N137").Value = -25496.7


$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5427.3335
$ws.Range("I80").Value = 4005
$ws.Range("K80").Value = 4005
$ws.Range("M80").Value = -3007

# Row 83
$ws.Range("H83").Value = 5427.3335
$ws.Range("I83").Value = 4005
$ws.Range("K83").Value = 20025
$ws.Range("M83").Value = -15033

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
# Row 24
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

# Row 68
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

# Row 93
$ws.Range("H93").Value = 111113656
$ws.Range("I93").Value = 142859220
$ws.Range("J93").Value = 4197.5
$ws.Range("K93").Value = 142859220
$ws.Range("L93").Value = 4197.5
$ws.Range("M93").Value = -142857972
$ws.Range("N93").Value = -6693.5


$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 2491666
$ws.Range("I8").Value = 999999
$ws.Range("J8").Value = 3983333
$ws.Range("K8").Value = 999999
$ws.Range("L8").Value = 3983333
$ws.Range("M8").Value = -999859
$ws.Range("N8").Value = -3983613

# Row 20
$ws.Range("H20").Value = 7500
$ws.Range("J20").Value = 7500
$ws.Range("L20").Value = 7500
$ws.Range("N20").Value = -7980

# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()

# Row 113
$ws.Range("H113").Value = 1465.75
$ws.Range("I113").Value = 1313.4286
$ws.Range("J113").Value = 1679
$ws.Range("K113").Value = 3940.2858
$ws.Range("L113").Value = 5037
$ws.Range("M113").Value = -1770.2858
$ws.Range("N113").Value = -9377

# Row 126
$ws.Range("H126").Value = 2642.5
$ws.Range("I126").Value = 2642.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7927.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5457.5
$ws.Range("N126").ClearContents()

# Row 136
$ws.Range("H136").Value = 3000.3333
$ws.Range("I136").Value = 2961.303
$ws.Range("K136").Value = 8883.909
$ws.Range("M136").Value = -6333.909

